$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 105.25
$ws.Range("I4").Value = 105.25
$ws.Range("K4").Value = 105.25
$ws.Range("M4").Value = 8.75
$ws.Range("H64").Value = 3300
$ws.Range("J64").Value = 3300
$ws.Range("L64").Value = 3300
$ws.Range("N64").Value = -3796
$ws.Range("H67").Value = 3300
$ws.Range("J67").Value = 3300
$ws.Range("L67").Value = 3300
$ws.Range("N67").Value = -5016
$ws.Range("H129").Value = 124387.38
$ws.Range("J129").Value = 137990.45
$ws.Range("L129").Value = 413971.35
$ws.Range("N129").Value = -423971.35

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2268
$ws.Range("I61").Value = 1907.7894
$ws.Range("J61").Value = 3123.5
$ws.Range("K61").Value = 1907.7894
$ws.Range("L61").Value = 3123.5
$ws.Range("M61").Value = -1695.7894
$ws.Range("N61").Value = -3547.5
$ws.Range("H74").Value = 27779192
$ws.Range("I74").Value = 34483364
$ws.Range("K74").Value = 34483364
$ws.Range("M74").Value = -34482490
$ws.Range("H77").Value = 27779192
$ws.Range("I77").Value = 34483364
$ws.Range("K77").Value = 172416820
$ws.Range("M77").Value = -172412452
$ws.Range("H122").Value = 1908.05
$ws.Range("I122").Value = 1734.7894
$ws.Range("K122").Value = 5204.3682
$ws.Range("M122").Value = -2754.3682
$ws.Range("H132").Value = 15491.564
$ws.Range("I132").Value = 2564.6667
$ws.Range("J132").Value = 44577.082
$ws.Range("K132").Value = 7694.000100000001
$ws.Range("L132").Value = 133731.246
$ws.Range("M132").Value = -5164.000100000001
$ws.Range("N132").Value = -138791.246
$ws.Range("H133").Value = 71393
$ws.Range("J133").Value = 71393
$ws.Range("L133").Value = 71393
$ws.Range("N133").Value = -76453
$ws.Range("H136").Value = 2268
$ws.Range("I136").Value = 1907.7894
$ws.Range("J136").Value = 3123.5
$ws.Range("K136").Value = 5723.3682
$ws.Range("L136").Value = 9370.5
$ws.Range("M136").Value = -3173.3682
$ws.Range("N136").Value = -14470.5
$ws.Range("H139").Value = 40692.25
$ws.Range("J139").Value = 40692.25
$ws.Range("L139").Value = 40692.25
$ws.Range("N139").Value = -50972.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 35829.4
$ws.Range("J82").Value = 42722.5
$ws.Range("L82").Value = 42722.5
$ws.Range("N82").Value = -43488.5
$ws.Range("H85").Value = 35829.4
$ws.Range("J85").Value = 42722.5
$ws.Range("L85").Value = 42722.5
$ws.Range("N85").Value = -45374.5
$ws.Range("H86").Value = 1840.909
$ws.Range("I86").Value = 1700.1818
$ws.Range("K86").Value = 1700.1818
$ws.Range("M86").Value = -577.1818000000001
$ws.Range("H89").Value = 1840.909
$ws.Range("I89").Value = 1700.1818
$ws.Range("K89").Value = 8500.909
$ws.Range("M89").Value = -2884.909
$ws.Range("H105").Value = 894664.25
$ws.Range("I105").Value = 1564.9231
$ws.Range("K105").Value = 1564.9231
$ws.Range("M105").Value = 182.0769
$ws.Range("H134").Value = 4444.5
$ws.Range("I134").Value = 5229.3687
$ws.Range("J134").Value = 1462
$ws.Range("K134").Value = 15688.1061
$ws.Range("L134").Value = 4386
$ws.Range("M134").Value = -13153.1061
$ws.Range("N134").Value = -9456

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3870.186
$ws.Range("I31").Value = 2652.5
$ws.Range("J31").Value = 4458.0347
$ws.Range("K31").Value = 2652.5
$ws.Range("L31").Value = 4458.0347
$ws.Range("M31").Value = -2357.5
$ws.Range("N31").Value = -5048.0347
$ws.Range("H34").Value = 3870.186
$ws.Range("I34").Value = 2652.5
$ws.Range("J34").Value = 4458.0347
$ws.Range("K34").Value = 2652.5
$ws.Range("L34").Value = 4458.0347
$ws.Range("M34").Value = -2450.5
$ws.Range("N34").Value = -4862.0347
$ws.Range("H138").Value = 34991.11
$ws.Range("J138").Value = 34991.11
$ws.Range("L138").Value = 34991.11
$ws.Range("N138").Value = -45271.11

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1205.9
$ws.Range("I18").Value = 211.8
$ws.Range("J18").Value = 2200
$ws.Range("K18").Value = 635.4000000000001
$ws.Range("L18").Value = 6600
$ws.Range("M18").Value = -466.4000000000001
$ws.Range("N18").Value = -6938
$ws.Range("H56").Value = 6506.923
$ws.Range("I56").Value = 6506.923
$ws.Range("K56").Value = 6506.923
$ws.Range("M56").Value = -5976.923
$ws.Range("H68").Value = 21080.4
$ws.Range("J68").Value = 26100.5
$ws.Range("L68").Value = 78301.5
$ws.Range("N68").Value = -79923.5
$ws.Range("H71").Value = 21080.4
$ws.Range("J71").Value = 26100.5
$ws.Range("L71").Value = 234904.5
$ws.Range("N71").Value = -243016.5
$ws.Range("H121").Value = 5051429
$ws.Range("I121").Value = 300
$ws.Range("J121").Value = 5209276.5
$ws.Range("K121").Value = 900
$ws.Range("L121").Value = 15627829.5
$ws.Range("M121").Value = 410
$ws.Range("N121").Value = -15630449.5
$ws.Range("H131").Value = 726.67
$ws.Range("J131").Value = 732.95874
$ws.Range("L131").Value = 2198.87622
$ws.Range("N131").Value = -12278.87622

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4670.077
$ws.Range("I113").Value = 5386.65
$ws.Range("J113").Value = 2281.5
$ws.Range("K113").Value = 5386.65
$ws.Range("L113").Value = 2281.5
$ws.Range("M113").Value = -3216.65
$ws.Range("N113").Value = -6621.5
$ws.Range("H122").Value = 3401.6316
$ws.Range("I122").Value = 2437.3572
$ws.Range("J122").Value = 6101.6
$ws.Range("K122").Value = 7312.071599999999
$ws.Range("L122").Value = 18304.8
$ws.Range("M122").Value = -4862.071599999999
$ws.Range("N122").Value = -23204.8
$ws.Range("H126").Value = 4847.8213
$ws.Range("I126").Value = 5600.9165
$ws.Range("J126").Value = 4283
$ws.Range("K126").Value = 16802.7495
$ws.Range("L126").Value = 12849
$ws.Range("M126").Value = -14332.7495
$ws.Range("N126").Value = -17789
$ws.Range("H132").Value = 23192.68
$ws.Range("I132").Value = 2681.8
$ws.Range("K132").Value = 8045.400000000001
$ws.Range("M132").Value = -5515.400000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4277.409
$ws.Range("I7").Value = 4257
$ws.Range("J7").Value = 4313.125
$ws.Range("K7").Value = 4257
$ws.Range("L7").Value = 4313.125
$ws.Range("M7").Value = -4145
$ws.Range("N7").Value = -4537.125
$ws.Range("H35").Value = 1000031
$ws.Range("I35").Value = 1000031
$ws.Range("K35").Value = 1000031
$ws.Range("M35").Value = -999695
$ws.Range("H61").Value = 3152.5715
$ws.Range("I61").Value = 1719.3529
$ws.Range("K61").Value = 1719.3529
$ws.Range("M61").Value = -1517.3529
$ws.Range("H100").Value = 1498.4
$ws.Range("I100").Value = 897.3333
$ws.Range("J100").Value = 2400
$ws.Range("K100").Value = 897.3333
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -356.3333
$ws.Range("N100").Value = -3482
$ws.Range("H111").Value = 29137
$ws.Range("J111").Value = 29137
$ws.Range("L111").Value = 29137
$ws.Range("N111").Value = -37317
$ws.Range("H113").Value = 3152.5715
$ws.Range("I113").Value = 1719.3529
$ws.Range("K113").Value = 1719.3529
$ws.Range("M113").Value = 450.6470999999999
$ws.Range("H126").Value = 4277.409
$ws.Range("I126").Value = 4257
$ws.Range("J126").Value = 4313.125
$ws.Range("K126").Value = 12771
$ws.Range("L126").Value = 12939.375
$ws.Range("M126").Value = -10301
$ws.Range("N126").Value = -17879.375
$ws.Range("H132").Value = 418067.4
$ws.Range("I132").Value = 525519.9399999999
$ws.Range("J132").Value = 6166.1665
$ws.Range("K132").Value = 1576559.82
$ws.Range("L132").Value = 18498.4995
$ws.Range("M132").Value = -1574029.82
$ws.Range("N132").Value = -23558.4995
$ws.Range("H136").Value = 1459.48
$ws.Range("I136").Value = 1570.091
$ws.Range("J136").Value = 648.3333
$ws.Range("K136").Value = 4710.272999999999
$ws.Range("L136").Value = 1944.9999
$ws.Range("M136").Value = -2160.272999999999
$ws.Range("N136").Value = -7044.9999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1934.2
$ws.Range("I81").Value = 930.3333
$ws.Range("J81").Value = 3440
$ws.Range("K81").Value = 1860.6666
$ws.Range("L81").Value = 6880
$ws.Range("M81").Value = -799.6666
$ws.Range("N81").Value = -9002
$ws.Range("H84").Value = 1934.2
$ws.Range("I84").Value = 930.3333
$ws.Range("J84").Value = 3440
$ws.Range("K84").Value = 9303.333000000001
$ws.Range("L84").Value = 34400
$ws.Range("M84").Value = -3999.333000000001
$ws.Range("N84").Value = -45008
$ws.Range("H126").Value = 1831.3438
$ws.Range("I126").Value = 1311.6538
$ws.Range("K126").Value = 3934.9614
$ws.Range("M126").Value = -1464.9614
$ws.Range("H132").Value = 1280.5217
$ws.Range("I132").Value = 738.46155
$ws.Range("K132").Value = 2215.38465
$ws.Range("M132").Value = 314.61535
$ws.Range("H136").Value = 30363912
$ws.Range("I136").Value = 38233204
$ws.Range("J136").Value = 10929.286
$ws.Range("K136").Value = 114699612
$ws.Range("L136").Value = 32787.858
$ws.Range("M136").Value = -114697062
